# Update the "Recorded By" column (G) values for several rows:
# the list of recorders was reordered in the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Amira.Sobhy@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G3").Value = "System, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G4").Value = "servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G7").Value = "Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg"
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"

$wb.Save()
